$wb = $excel.ActiveWorkbook

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 2663.2144
$ws.Range("I118").Value = 581.6667
$ws.Range("J118").Value = 4224.375
$ws.Range("K118").Value = 1745.0001
$ws.Range("L118").Value = 12673.125
$ws.Range("M118").Value = -88.00009999999997
$ws.Range("N118").Value = -15987.125

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2819101.2
$ws.Range("I132").Value = 3280713.2
$ws.Range("K132").Value = 9842139.600000001
$ws.Range("M132").Value = -9839609.600000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5252.9106
$ws.Range("I138").Value = 1499.8518
$ws.Range("J138").Value = 8747.138000000001
$ws.Range("K138").Value = 4499.555399999999
$ws.Range("L138").Value = 26241.414
$ws.Range("M138").Value = 640.4446000000007
$ws.Range("N138").Value = -36521.414

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 66
$ws.Range("N4").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2765.973
$ws.Range("I132").Value = 2204.25
$ws.Range("K132").Value = 6612.75
$ws.Range("M132").Value = -4082.75

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2314.2964
$ws.Range("I99").Value = 1700.3
$ws.Range("K99").Value = 1700.3
$ws.Range("M99").Value = -202.3

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3554.5264
$ws.Range("I31").Value = 2050.92
$ws.Range("K31").Value = 2050.92
$ws.Range("M31").Value = -1755.92

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3554.5264
$ws.Range("I34").Value = 2050.92
$ws.Range("K34").Value = 2050.92
$ws.Range("M34").Value = -1848.92

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6581121
$ws.Range("I58").Value = 1450.5416
$ws.Range("J58").Value = 17860556
$ws.Range("K58").Value = 1450.5416
$ws.Range("L58").Value = 17860556
$ws.Range("M58").Value = -1247.5416
$ws.Range("N58").Value = -17860962

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 18797.285
$ws.Range("J59").Value = 18797.285
$ws.Range("L59").Value = 18797.285
$ws.Range("N59").Value = -21087.285

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1655.661
$ws.Range("I132").Value = 1202.4
$ws.Range("K132").Value = 3607.2
$ws.Range("M132").Value = -1077.2

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1695.6666
$ws.Range("I134").Value = 1086.2142
$ws.Range("K134").Value = 3258.6426
$ws.Range("M134").Value = -723.6425999999997

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6581121
$ws.Range("I136").Value = 1450.5416
$ws.Range("J136").Value = 17860556
$ws.Range("K136").Value = 4351.6248
$ws.Range("L136").Value = 53581668
$ws.Range("M136").Value = -1801.6248
$ws.Range("N136").Value = -53586768

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 25142.857
$ws.Range("J141").Value = 25142.857
$ws.Range("L141").Value = 25142.857
$ws.Range("N141").Value = -35502.857

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2426.111
$ws.Range("I131").Value = 3518.5715
$ws.Range("J131").Value = 1730.909
$ws.Range("K131").Value = 10555.7145
$ws.Range("L131").Value = 5192.727000000001
$ws.Range("M131").Value = -5515.7145
$ws.Range("N131").Value = -15272.727

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 307504.5
$ws.Range("J18").Value = 76672.664
$ws.Range("L18").Value = 76672.664
$ws.Range("N18").Value = -77258.664

# GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1957.6
$ws.Range("I43").Value = 1957.6
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1957.6
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1806.6
$ws.Range("N43").ClearContents()

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9000
$ws.Range("I57").Value = 9000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 9000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -8180
$ws.Range("N57").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2799.9
$ws.Range("I80").Value = 2799.9
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2799.9
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1801.9
$ws.Range("N80").ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2799.9
$ws.Range("I83").Value = 2799.9
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13999.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -9007.5
$ws.Range("N83").ClearContents()

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1586.3334
$ws.Range("I113").Value = 1035.5385
$ws.Range("J113").Value = 5166.5
$ws.Range("K113").Value = 1035.5385
$ws.Range("L113").Value = 5166.5
$ws.Range("M113").Value = 1134.4615
$ws.Range("N113").Value = -9506.5

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3031.7896
$ws.Range("I126").Value = 878
$ws.Range("J126").Value = 3801
$ws.Range("K126").Value = 2634
$ws.Range("L126").Value = 11403
$ws.Range("M126").Value = -164
$ws.Range("N126").Value = -16343

# LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 25232.615
$ws.Range("J20").Value = 25232.615
$ws.Range("L20").Value = 25232.615
$ws.Range("N20").Value = -25684.615

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1471.5834
$ws.Range("I46").Value = 554
$ws.Range("J46").Value = 1849.4117
$ws.Range("K46").Value = 554
$ws.Range("L46").Value = 1849.4117
$ws.Range("M46").Value = -366
$ws.Range("N46").Value = -2225.4117

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2022.52
$ws.Range("I68").Value = 1034.5454
$ws.Range("J68").Value = 9267.666999999999
$ws.Range("K68").Value = 1034.5454
$ws.Range("L68").Value = 9267.666999999999
$ws.Range("M68").Value = -285.5454
$ws.Range("N68").Value = -10765.667

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2022.52
$ws.Range("I71").Value = 1034.5454
$ws.Range("J71").Value = 9267.666999999999
$ws.Range("K71").Value = 5172.727
$ws.Range("L71").Value = 46338.335
$ws.Range("M71").Value = -1428.727
$ws.Range("N71").Value = -53826.335

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1798.6154
$ws.Range("I132").Value = 1131.7755
$ws.Range("K132").Value = 3395.3265
$ws.Range("M132").Value = -865.3265000000001

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 932.3333
$ws.Range("I107").Value = 492.6
$ws.Range("J107").Value = 2031.6666
$ws.Range("K107").Value = 1477.8
$ws.Range("L107").Value = 6094.9998
$ws.Range("M107").Value = 442.1999999999998
$ws.Range("N107").Value = -9934.9998

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1343.24
$ws.Range("I113").Value = 650.75
$ws.Range("J113").Value = 2574.3333
$ws.Range("K113").Value = 1952.25
$ws.Range("L113").Value = 7722.999899999999
$ws.Range("M113").Value = 217.75
$ws.Range("N113").Value = -12062.9999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12871.02
$ws.Range("I132").Value = 3216.6667
$ws.Range("K132").Value = 9650.000100000001
$ws.Range("M132").Value = -7120.000100000001
